$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.142.65"
$ws.Range("E2").Value = "  -4.37%  "
$ws.Range("D3").Value = "1.651.53"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'215.39"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").Value = "'0.5125"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.2592"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "'0.06434"
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("D10").Value = "'19.91"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "'0.07766"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "1.656.34"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "'4.285"
$ws.Range("E13").Value = "  -4.87%  "
$ws.Range("D14").Value = "1.879.41"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "'0.5515"
$ws.Range("E15").Value = "  -5.78%  "
$ws.Range("D16").Value = "0.0₅8005"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "'64.04"
$ws.Range("E17").Value = "  -5.73%  "
$ws.Range("D18").Value = "26.150.09"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'210.34"
$ws.Range("E20").Value = "  -5.58%  "
$ws.Range("D21").Value = "'4.389"
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").Value = "'10.05"
$ws.Range("E22").Value = "  -4.08%  "
$ws.Range("D23").Value = "'6.038"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'1.744"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").Value = "'0.1174"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'6.969"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").Value = "'15.80"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "'0.05129"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "'3.350"
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").Value = "'3.211"
$ws.Range("E33").Value = "  -6.27%  "
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D35").Value = "'2.739"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.354"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9227"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").Value = "1.168.14"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").Value = "'0.5685"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "'0.01583"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'2.548"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'5.658"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "'0.8231"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").Value = "'99.93"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "1.789.44"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").Value = "0.0₈116"
$ws.Range("E47").Value = "  +3.71%  "
$ws.Range("D48").Value = "'0.4557"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'55.45"
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D51").Value = "'7.840"
$ws.Range("E51").Value = "  -3.19%  "
